$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-shuffles the weekly price records across rows 2-9 (columns
# D, L, M, N, O, P, Q, S). All other columns (A,B,C,E,F,G,H,I,J,K,R,T) are
# identical for every row and remain unchanged.

$rows = @{
    2 = @{ D = 44699; L = "Primera"; M = 100; N = 20000; O = 22000; P = 21000; Q = "`$/caja 18 kilos";        S = 1167 }
    3 = @{ D = 44699; L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos";        S = 1000 }
    4 = @{ D = 44516; L = "Primera"; M = 100; N = 33000; O = 34000; P = 33500; Q = "`$/caja 18 kilos";        S = 1861 }
    5 = @{ D = 44316; L = "Primera"; M = 50;  N = 20000; O = 20000; P = 20000; Q = "`$/caja 18 kilos";        S = 1111 }
    6 = @{ D = 44819; L = "Primera"; M = 100; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos granel"; S = 1417 }
    7 = @{ D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos";        S = 806  }
    8 = @{ D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";        S = 667  }
    9 = @{ D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "`$/caja 18 kilos";        S = 1028 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("S$r").Value = $vals.S
}
